# "Add files via upload" — adds a "rating" column and a "products" column
# to the "productos" sheet, and removes a stray formatted-but-empty cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert the new "rating" column (new column E, pushing contact..oferta right) ---
$ws.Columns("E").Insert()

# --- Insert the new "products" column (new column I, pushing oferta..señora right) ---
$ws.Columns("I").Insert()

# The freshly-inserted column inherited its left neighbour's format; the header
# cell should instead carry the bold/bordered header style used by J1:N1.
$ws.Range("J1").Copy()
$ws.Range("I1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Column E: rating ---
$ws.Range("E1").Value = "rating"
$ws.Range("E2").Value = "4.6"
$ws.Range("E3").Value = "4.2"
$ws.Range("E4").Value = "3.5"
$ws.Range("E5").Value = "2.5"
$ws.Range("E6").Value = "4.5"
$ws.Range("E7").Value = "4.6"
$ws.Range("E8").Value = "3.6"
$ws.Range("E9").Value = "4.5"
$ws.Range("E10").Value = "4.6"
$ws.Range("E11").Value = "3.6"

# --- Column I: products (only the first four product rows got data) ---
$ws.Range("I1").Value = "products"
$ws.Range("I2").Value = "huevos, pizza, arroz"
$ws.Range("I3").Value = "pasta, arroz, frutas"
$ws.Range("I4").Value = "mango, dulces, pan"
$ws.Range("I5").Value = "arina, trigo"

# Rows 6-11 never had a "products" entry; the column insert left behind blank
# placeholder cells there (copied formatting, no content) — drop them so the
# row only spans as far as it actually has data, same as the source rows.
$ws.Range("I6:I11").Clear()

# --- Remove the stray formatted cell at C15 (style-only cell, no value) ---
$ws.Range("C15").Clear()

# --- Restore the view/selection state recorded in the saved workbook ---
$win = $excel.ActiveWindow
$win.ScrollColumn = 6
$win.ScrollRow = 1
$ws.Range("J13").Select()
